$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.528.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = "'1.728.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'244.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.36%  '
$ws.Range("D7").Value = "'0.4795"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.76%  '
$ws.Range("D8").Value = "'0.2670"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.48%  '
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").Value = "'1.731.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.88%  '
$ws.Range("D11").Value = "'0.07161"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.58%  '
$ws.Range("D12").Value = "'15.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.23%  '
$ws.Range("D13").Value = "'0.6141"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.17%  '
$ws.Range("D14").Value = "'4.525"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.44%  '
$ws.Range("D15").Value = "'76.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.13%  '
$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").Value = "'26.531.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("D18").Value = "'1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = "'0.000006969"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.14%  '
$ws.Range("D20").Value = "'11.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.65%  '
$ws.Range("D21").Value = "'1.953.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.94%  '
$ws.Range("E22").Value = '  -0.49%  '
$ws.Range("D23").Value = "'8.904"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.69%  '
$ws.Range("D24").Value = "'5.283"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.69%  '
$ws.Range("D25").Value = "'136.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.71%  '
$ws.Range("D26").Value = "'15.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.54%  '
$ws.Range("D27").Value = "'1.790"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.17%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").Value = "'106.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.34%  '
$ws.Range("D30").Value = "'3.979"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").Value = "'0.07963"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.87%  '
$ws.Range("D32").Value = "'3.707"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.49%  '
$ws.Range("D33").Value = "'0.04577"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.06%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("D36").Value = "'0.9940"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.00%  '
$ws.Range("D37").Value = "'0.6316"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.00%  '
$ws.Range("D38").Value = "'2.089"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.79%  '
$ws.Range("D39").Value = "'0.9148"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.05%  '
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("D41").Value = "'104.65"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.64%  '
$ws.Range("D42").Value = "'1.003"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("E43").Value = '  +2.01%  '
$ws.Range("D44").Value = "'5.574"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.49%  '
$ws.Range("E45").Value = '  +1.73%  '
$ws.Range("D46").Value = "'6.992"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.12%  '
$ws.Range("E47").Value = '  +1.93%  '
$ws.Range("D48").Value = "'0.05345"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.16%  '
$ws.Range("D49").Value = "'30.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.24%  '
$ws.Range("D50").Value = "'7.863"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.49%  '
$ws.Range("D51").Value = "'1.258"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.71%  '
